$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '98.752.34'
$ws.Range('E2').Value = '  -0.37%  '
$ws.Range('D3').Value = '3.336.88'
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '255.53'
$ws.Range('E5').Value = '  -2.69%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '642.20'
$ws.Range('E6').Value = '  +1.14%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.57'
$ws.Range('E7').Value = '  +14.37%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.424'
$ws.Range('E8').Value = '  +7.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.09'
$ws.Range('E9').Value = '  +25.13%  '
$ws.Range('E10').Value = '  -0.05%  '
$ws.Range('D11').Value = '3.336.46'
$ws.Range('E11').Value = '  -0.25%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.207'
$ws.Range('E12').Value = '  +3.02%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '43.49'
$ws.Range('E13').Value = '  +18.90%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000272'
$ws.Range('E14').Value = '  +8.77%  '
$ws.Range('D15').Value = '98.572.78'
$ws.Range('E15').Value = '  -0.05%  '
$ws.Range('E16').Value = '  -0.30%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.52'
$ws.Range('E17').Value = '  -0.60%  '
$ws.Range('D18').Value = '3.316.15'
$ws.Range('E18').Value = '  -0.43%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.04'
$ws.Range('E19').Value = '  +14.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.67'
$ws.Range('E20').Value = '  +10.94%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '543.58'
$ws.Range('E21').Value = '  +9.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.49'
$ws.Range('E22').Value = '  -2.31%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.13'
$ws.Range('E23').Value = '  +8.90%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.440'
$ws.Range('E24').Value = '  +57.94%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000200'
$ws.Range('E25').Value = '  -4.82%  '
$ws.Range('B26').Value = 'Litecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '100.58'
$ws.Range('E26').Value = '  +13.26%  '
$ws.Range('B27').Value = 'NEARProtocol'
$ws.Range('C27').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.11'
$ws.Range('E27').Value = '  +6.56%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.51'
$ws.Range('E28').Value = '  +3.89%  '
$ws.Range('D29').Value = '3.506.09'
$ws.Range('E29').Value = '  -0.65%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.152'
$ws.Range('E30').Value = '  +19.22%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.999'
$ws.Range('E31').Value = '  -0.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '10.97'
$ws.Range('E32').Value = '  +15.19%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.190'
$ws.Range('E33').Value = '  -3.49%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '29.26'
$ws.Range('E35').Value = '  +5.41%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.523'
$ws.Range('E36').Value = '  +13.10%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '7.59'
$ws.Range('E37').Value = '  +2.76%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.155'
$ws.Range('E38').Value = '  +3.99%  '
$ws.Range('B39').Value = 'PancakeSwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.04'
$ws.Range('E39').Value = '  +2.50%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '520.14'
$ws.Range('E40').Value = '  +2.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '24.73'
$ws.Range('E41').Value = '  -0.50%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.85'
$ws.Range('E42').Value = '  +3.24%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.30'
$ws.Range('E43').Value = '  +1.61%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.810'
$ws.Range('E44').Value = '  +5.77%  '
$ws.Range('E45').Value = '  -0.02%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0393'
$ws.Range('E46').Value = '  +22.32%  '
$ws.Range('B47').Value = 'dogwifhat'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.17'
$ws.Range('E47').Value = '  -2.80%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.03'
$ws.Range('E48').Value = '  +3.95%  '
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.76'
$ws.Range('E49').Value = '  +18.31%  '
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '163.51'
$ws.Range('E50').Value = '  +1.26%  '
$ws.Range('B51').Value = 'OKB'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '49.40'
$ws.Range('E51').Value = '  +5.34%  '
